$wb = $excel.ActiveWorkbook

# --- Sheet: Inventory (Stock column D) ---
$inv = $wb.Worksheets.Item("Inventory")
$inv.Range("D2").Value = 991
$inv.Range("D3").Value = 793
$inv.Range("D5").Value = 820
$inv.Range("D6").Value = 9982
$inv.Range("D7").Value = 9820
$inv.Range("D9").Value = 55
$inv.Range("D10").Value = 91
$inv.Range("D11").Value = 0
$inv.Range("F12").Select()

# --- Sheet: Shopping List (add row 12 for Joe / Laptop) ---
$sl = $wb.Worksheets.Item("Shopping List")
$sl.Range("A12").Value = "Joe"
$sl.Range("B12").Value = "Laptop"
$sl.Range("C12").Value = 2
$sl.Range("D14").Select()

# --- Sheet: Expenses (update Joe's total) ---
$exp = $wb.Worksheets.Item("Expenses")
$exp.Range("B4").Value = 758.75
$exp.Range("A2:B4").Select()

# --- Sheet: Items Not Found (add row 3 for Joe / Laptop) ---
$inf = $wb.Worksheets.Item("Items Not Found")
$inf.Range("A3").Value = "Joe"
$inf.Range("B3").Value = "Laptop"
$inf.Range("C3").Value = 1
$inf.Range("E18").Select()

# Restore Expenses as the active/selected tab (it was the active sheet originally)
$exp.Select()
$exp.Range("A2:B4").Select()

$wb.Save()
